$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5
$ws.Range("A5").Value = "Wizard"
$ws.Range("B5").Value = "Tome Master"
$ws.Range("E5").Value = "ALL"
$ws.Range("J5").Value = "Fringe Mage"
$ws.Range("O5").Value = "ALL"
$ws.Range("U5").Value = "ALL"
$ws.Range("V5").Value = "Elementalist"

# Row 6
$ws.Range("A6").Value = "Marksman"
$ws.Range("C6").Value = "ALL"
$ws.Range("D6").Value = "Zen Archer"
$ws.Range("F6").Value = "Juggernaut"
$ws.Range("H6").Value = "ALL"
$ws.Range("L6").Value = "Doom Sniper"
$ws.Range("P6").Value = "ALL"

# Row 7
$ws.Range("A7").Value = "Pugilist"
$ws.Range("D7").Value = "ALL"
$ws.Range("J7").Value = "Street Brawler"
$ws.Range("M7").Value = "ALL"
$ws.Range("N7").Value = "ALL"
$ws.Range("Q7").Value = "Monk"
$ws.Range("V7").Value = "Fire Fist"

# Column J width adjustment (target OOXML width 12.96 chars; engine quantizes
# ColumnWidth to a whole-pixel grid before re-deriving the stored width, so
# 12.17 is the nearest settable value that lands on the same pixel column as
# 12.96 would)
$ws.Columns.Item(10).ColumnWidth = 12.17

# Selection moves to A8
$ws.Range("A8").Select()
